# Rechtschreibfehler + Regelbetrieb Definition
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "Papierkneuel" -> "Papierknäuel"
$ws.Range("A2").Value = "Papierknäuel"

# Fix Regelbetrieb definition wording
$ws.Range("B17").Value = "½ Stunde ununterbrochen in einen Job tätig."

# Fix "Piken" -> "Pieken" and clarify the definition text
$ws.Range("A19").Value = "Pieken"
$ws.Range("B19").Value = "Unmittelbar aufeinanderfolgende Ab- und Aufwärtsbewegung zur Aufnahme von Zigarettenstummel"
